$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.699297666549683
$ws.Range("B1").Value = 3.186416387557983
$ws.Range("C1").Value = 2.468589305877686
$ws.Range("D1").Value = 2.288189172744751
$ws.Range("E1").Value = 1.922728180885315
